$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.00793993473053
$ws.Range("B1").Value = 1.813161492347717
$ws.Range("C1").Value = 5.106608390808105
$ws.Range("D1").Value = 1.484330773353577
$ws.Range("E1").Value = 1.360289931297302
